# Apply the edits described by the commit "update AR excel files and html files"
# to the Sheet1 of bubble_discretionary.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a handful of cells that were showing raw numbers instead of the
#     intended percentage format (0.00%) -------------------------------
$ws.Range("B9").NumberFormat = "0.00%"
$ws.Range("B18").NumberFormat = "0.00%"

# E17 was stored as a raw number (71.5, formatted like currency) instead of
# the fraction (0.715) formatted as a percentage.
$ws.Range("E17").NumberFormat = "0.00%"
$ws.Range("E17").Value = 0.715

# --- Fill in the previously blank TOTALS row percentages ----------------
$ws.Range("B22:E22").Font.Size = 10
$ws.Range("B22").NumberFormat = "0.00%"
$ws.Range("C22").NumberFormat = "0.00%"
$ws.Range("E22").NumberFormat = "0.00%"
$ws.Range("B22").Value = 0.056
$ws.Range("C22").Value = 1
$ws.Range("E22").Value = 1

# --- Update the sheet's current selection (cosmetic, matches re-save) ---
$ws.Range("F24").Select()

# --- Set the page to print in portrait orientation -----------------------
$ws.PageSetup.Orientation = 1
